# Corrections following third round of review
#
# The "Materials" sheet carried a Darwin Core "subgenus" term between the
# "genus" and "specificEpithet" columns (column AS). This term is being
# removed from the template entirely: the whole column (header in row 1
# and the "${subgenus}" placeholder in row 2) is deleted and the columns
# to its right shift one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Column AS holds the "subgenus" header (row 1) / "${subgenus}" token (row 2).
# Deleting the entire column shifts everything after it one place left,
# matching the EZ/EY dimension collapse seen in the diff, and also drops
# the now-unused "subgenus" / "${subgenus}" shared strings automatically.
$ws.Range("AS1").EntireColumn.Delete()
